# Update countries & provincias Spain
# Refresh the "Pais" data: bump the "Datos actualizados" timestamp and
# rewrite the per-country statistics that moved as the case totals were
# refreshed (several countries overtook their neighbours in the
# case-count ranking, so their rows swapped places).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp footer
$ws.Range("A1").Value = "Datos actualizados a 28 de Mayo de 2020 a las 11:10"

# Alemania (row 11) - Casos activos / Recuperados refreshed
$ws.Range("D11").Value = 163200
$ws.Range("E11").Value = 10162

# Banglades overtakes Bielorrusia (rows 25-26 swap)
$ws.Range("A25").Value = "Banglades"
$ws.Range("B25").Value = 40321
$ws.Range("C25").Value = 2029
$ws.Range("D25").Value = 8425
$ws.Range("E25").Value = 31337
$ws.Range("G25").Value = 15
$ws.Range("H25").Value = 559

$ws.Range("A26").Value = "Bielorrusia"
$ws.Range("B26").Value = 38956
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 15923
$ws.Range("E26").Value = 22819
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = 214

# Indonesia overtakes Colombia (rows 35-36 swap)
$ws.Range("A35").Value = "Indonesia"
$ws.Range("B35").Value = 24538
$ws.Range("C35").Value = 687
$ws.Range("D35").Value = 6240
$ws.Range("E35").Value = 16802
$ws.Range("G35").Value = 23
$ws.Range("H35").Value = 1496

$ws.Range("A36").Value = "Colombia"
$ws.Range("B36").Value = 24104
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 6111
$ws.Range("E36").Value = 17190
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 803

# Finlandia (row 67) - Casos totales / Nuevos casos / Recuperados refreshed
$ws.Range("B67").Value = 6743
$ws.Range("C67").Value = 51
$ws.Range("E67").Value = 1330

# Nepal jumps ahead of Sudan del Sur / Costa Rica / Niger / Republica de Chipre
# (rows 115-119 cascade down one position each)
$ws.Range("A115").Value = "Nepal"
$ws.Range("B115").Value = 1042
$ws.Range("C115").Value = 156
$ws.Range("D115").Value = 187
$ws.Range("E115").Value = 850
$ws.Range("G115").Value = 1
$ws.Range("H115").Value = 5

$ws.Range("A116").Value = "Sudan del Sur"
$ws.Range("B116").Value = 994
$ws.Range("C116").Value = 0
$ws.Range("D116").Value = 6
$ws.Range("E116").Value = 978
$ws.Range("F116").Value = 0
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 10

$ws.Range("A117").Value = "Costa Rica"
$ws.Range("B117").Value = 984
$ws.Range("C117").Value = 0
$ws.Range("D117").Value = 639
$ws.Range("E117").Value = 335
$ws.Range("F117").Value = 0
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 10

$ws.Range("A118").Value = "Niger"
$ws.Range("B118").Value = 952
$ws.Range("C118").Value = 0
$ws.Range("D118").Value = 796
$ws.Range("E118").Value = 93
$ws.Range("F118").Value = 0
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 63

$ws.Range("A119").Value = "Republica de Chipre"
$ws.Range("B119").Value = 939
$ws.Range("C119").Value = 0
$ws.Range("D119").Value = 594
$ws.Range("E119").Value = 328
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 17

# Belice overtakes Nueva Caledonia (rows 200-201 swap)
$ws.Range("A200").Value = "Belice"
$ws.Range("D200").Value = 16
$ws.Range("H200").Value = 2

$ws.Range("A201").Value = "Nueva Caledonia"
$ws.Range("D201").Value = 18
$ws.Range("H201").Value = 0

# Papua Nueva Guinea overtakes Islas Virgenes Britanicas (rows 213-214 swap)
$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("D213").Value = 8
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Islas Virgenes Britanicas"
$ws.Range("D214").Value = 7
$ws.Range("H214").Value = 1

# San Bartolome overtakes Bonaire, San Eustaquio y Saba (rows 215-216 swap, no stat change)
$ws.Range("A215").Value = "San Bartolome"
$ws.Range("A216").Value = "Bonaire, San Eustaquio y Saba"
